$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet index 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 150
$ws1.Range("F3").Value = 357
$ws1.Range("F4").Value = 455
$ws1.Range("F5").Value = 1786
$ws1.Range("F6").Value = 92
$ws1.Range("F7").Value = 2241
$ws1.Range("F9").Value = 290
$ws1.Range("F11").Value = 5069
$ws1.Range("F12").Value = 373
$ws1.Range("F14").Value = 314
$ws1.Range("F17").Value = 204
$ws1.Range("F21").Value = 4095
$ws1.Range("F22").Value = 736
$ws1.Range("F23").Value = 732
$ws1.Range("G23").Value = 55
$ws1.Range("F25").Value = 26
$ws1.Range("F26").Value = 116
$ws1.Range("F30").Value = 99
$ws1.Range("F34").Value = 1044
$ws1.Range("F36").Value = 2655
$ws1.Range("F37").Value = 432
$ws1.Range("F38").Value = 44

# --- Sheet "全部类型" (sheet index 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 150
$ws4.Range("F3").Value = 357
$ws4.Range("F4").Value = 455
$ws4.Range("F5").Value = 1786
$ws4.Range("F6").Value = 92
$ws4.Range("F7").Value = 2241
$ws4.Range("F9").Value = 290
$ws4.Range("F11").Value = 5069
$ws4.Range("F12").Value = 373
$ws4.Range("F14").Value = 314
$ws4.Range("F17").Value = 204
$ws4.Range("F21").Value = 4095
$ws4.Range("F22").Value = 736
$ws4.Range("F23").Value = 732
$ws4.Range("G23").Value = 55
$ws4.Range("F25").Value = 26
$ws4.Range("F26").Value = 116
$ws4.Range("F30").Value = 99
$ws4.Range("F35").Value = 1044
$ws4.Range("F37").Value = 2655
$ws4.Range("F38").Value = 432
$ws4.Range("F39").Value = 44
